$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = "297.24"
$ws.Range("E2").Value = "1.76%"
$ws.Range("D2:E2").Style = "Normal"

# Row 3
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = "41.86"
$ws.Range("E3").Value = "3.73%"
$ws.Range("D3:E3").Style = "Normal"

# Row 4
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = "5.003"
$ws.Range("E4").Value = "-0.14%"
$ws.Range("D4:E4").Style = "Normal"

# Row 5
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07519"
$ws.Range("E5").Value = "2.76%"
$ws.Range("D5:E5").Style = "Normal"

# Row 6
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = "1.584"
$ws.Range("E6").Value = "4.00%"
$ws.Range("D6:E6").Style = "Normal"

# Row 7
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9263"
$ws.Range("E7").Value = "-0.01%"
$ws.Range("D7:E7").Style = "Normal"

# Row 8
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "0.97%"
$ws.Range("E8").Style = "Normal"

# Row 9
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = "0.1200"
$ws.Range("E9").Value = "-0.33%"
$ws.Range("D9:E9").Style = "Normal"

# Row 10
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1834"
$ws.Range("E10").Value = "5.48%"
$ws.Range("D10:E10").Style = "Normal"

# Row 11
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08872"
$ws.Range("E11").Value = "2.88%"
$ws.Range("D11:E11").Style = "Normal"

# Row 12
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = "0.04076"
$ws.Range("E12").Value = "-5.95%"
$ws.Range("D12:E12").Style = "Normal"

# Row 13
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = "0.1049"
$ws.Range("E13").Value = "-0.47%"
$ws.Range("D13:E13").Style = "Normal"

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "1.24%"
$ws.Range("E14").Style = "Normal"

# Row 15
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = "0.005966"
$ws.Range("E15").Value = "-0.54%"
$ws.Range("D15:E15").Style = "Normal"

# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "0.56%"
$ws.Range("E16").Style = "Normal"

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "4.382"
$ws.Range("D17").Style = "Normal"

# Row 18
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = "0.3313"
$ws.Range("E18").Value = "0.76%"
$ws.Range("D18:E18").Style = "Normal"

# Row 19
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = "8.122"
$ws.Range("E19").Value = "4.90%"
$ws.Range("D19:E19").Style = "Normal"

# Row 20
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "0.04%"
$ws.Range("E20").Style = "Normal"

# Row 21
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = "0.3102"
$ws.Range("E21").Value = "10.99%"
$ws.Range("D21:E21").Style = "Normal"

# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "4.09%"
$ws.Range("E22").Style = "Normal"

# Row 23
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = "0.001266"
$ws.Range("E23").Value = "0.37%"
$ws.Range("D23:E23").Style = "Normal"

# Row 24
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = "0.003905"
$ws.Range("E24").Value = "3.36%"
$ws.Range("D24:E24").Style = "Normal"

# Row 25
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-3.95%"
$ws.Range("E25").Style = "Normal"

# Row 38
$ws.Range("D38:E38").NumberFormat = "@"
$ws.Range("D38").Value = "0.02404"
$ws.Range("E38").Value = "5.08%"
$ws.Range("D38:E38").Style = "Normal"

# Row 39
$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("D39").Value = "0.05235"
$ws.Range("E39").Value = "5.07%"
$ws.Range("D39:E39").Style = "Normal"

# Row 40
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "17.69%"
$ws.Range("E40").Style = "Normal"

# Row 41
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007826"
$ws.Range("E41").Value = "1.82%"
$ws.Range("D41:E41").Style = "Normal"

# Row 42
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1325"
$ws.Range("E42").Value = "3.04%"
$ws.Range("D42:E42").Style = "Normal"

# Row 43
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.98%"
$ws.Range("E43").Style = "Normal"

# Row 44
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = "0.007816"
$ws.Range("E44").Value = "-1.19%"
$ws.Range("D44:E44").Style = "Normal"

# Row 45
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = "0.2959"
$ws.Range("E45").Value = "-6.88%"
$ws.Range("D45:E45").Style = "Normal"

# Row 46
$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006526"
$ws.Range("E46").Value = "3.36%"
$ws.Range("D46:E46").Style = "Normal"

# Row 47
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "-0.05%"
$ws.Range("E47").Style = "Normal"

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "53.16%"
$ws.Range("E48").Style = "Normal"

# Row 49
$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("D49").Value = "0.004204"
$ws.Range("E49").Value = "0.07%"
$ws.Range("D49:E49").Style = "Normal"

# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.05%"
$ws.Range("E50").Style = "Normal"

# Row 51
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "-0.05%"
$ws.Range("E51").Style = "Normal"
